$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values (rows 2-8) to reflect the corrected/edited data
$ws.Range("A2").Value = 3000
$ws.Range("B2").Value = 3000

$ws.Range("A3").Value = 1000
$ws.Range("B3").Value = 1000

$ws.Range("A6").Value = 2000
$ws.Range("B6").Value = 2000

$ws.Range("A7").Value = 2000
$ws.Range("B7").Value = 2000

$ws.Range("A8").Value = 1400
$ws.Range("B8").Value = 1400

# Remove the last two rows (9 and 10), which are no longer part of the list
$ws.Range("A9:B10").Delete() | Out-Null
